# Apply the "Adds proper trained model to the examples" edit:
#   1. Translate the four sheet names to Russian.
#   2. Add row labels to the "Learning Factor" sheet, pushing the existing
#      numeric values from column A into column B.

$wb = $excel.ActiveWorkbook

# --- 1. Rename sheets -------------------------------------------------
$wsLearningFactor = $wb.Worksheets.Item("Learning Factor")
$wsTier1 = $wb.Worksheets.Item("Tier 1")
$wsTier2 = $wb.Worksheets.Item("Tier 2")
$wsTier3 = $wb.Worksheets.Item("Tier 3")

$wsLearningFactor.Name = "Коэффициент скорости обучения"
$wsTier1.Name = "Слой 1"
$wsTier2.Name = "Слой 2"
$wsTier3.Name = "Слой 3"

# --- 2. Add labels to the learning-factor sheet ------------------------
# Move the existing values from column A into column B (A1->B1, A2->B2),
# introduce a new B3 value, then label column A with row descriptions.
$wsLearningFactor.Range("B1").Value = 0.5
$wsLearningFactor.Range("B2").Value = 8.0
$wsLearningFactor.Range("B3").Value = 1.0

$wsLearningFactor.Range("A1").Value = "Коэффициент скорости обучения"
$wsLearningFactor.Range("A2").Value = "Количество входов"
$wsLearningFactor.Range("A3").Value = "Размерность выходного слоя"
